# Updates NATMI LR-pair TPM-derived metrics (Lrpap1-Lrp1) for rows 2-26,
# refreshing ligand/receptor average & total expression, their derived
# specificities, and the resulting edge weights/specificities.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 8.080435
$ws.Range("H2").Value = 24.241305
$ws.Range("I2").Value = 0.1496988574979475
$ws.Range("J2").Value = 0.1496988574979476
$ws.Range("M2").Value = 6.305846
$ws.Range("N2").Value = 18.917538
$ws.Range("O2").Value = 0.01356150511917599
$ws.Range("P2").Value = 0.01356150511917599
$ws.Range("Q2").Value = 50.95397872301
$ws.Range("R2").Value = 458.58580850709
$ws.Range("S2").Value = 0.002030141822293212
$ws.Range("T2").Value = 0.002030141822293212
$ws.Range("G3").Value = 8.080435
$ws.Range("H3").Value = 24.241305
$ws.Range("I3").Value = 0.1496988574979475
$ws.Range("J3").Value = 0.1496988574979476
$ws.Range("O3").Value = 0.392557056479861
$ws.Range("P3").Value = 0.3925570564798609
$ws.Range("Q3").Value = 1474.935394535148
$ws.Range("R3").Value = 13274.41855081634
$ws.Range("S3").Value = 0.05876534285779245
$ws.Range("T3").Value = 0.05876534285779245
$ws.Range("G4").Value = 8.080435
$ws.Range("H4").Value = 24.241305
$ws.Range("I4").Value = 0.1496988574979475
$ws.Range("J4").Value = 0.1496988574979476
$ws.Range("M4").Value = 127.396393
$ws.Range("N4").Value = 382.189179
$ws.Range("O4").Value = 0.2739817680029065
$ws.Range("P4").Value = 0.2739817680029065
$ws.Range("Q4").Value = 1029.418272870955
$ws.Range("R4").Value = 9264.764455838596
$ws.Range("S4").Value = 0.04101475764530282
$ws.Range("T4").Value = 0.04101475764530283
$ws.Range("G5").Value = 8.080435
$ws.Range("H5").Value = 24.241305
$ws.Range("I5").Value = 0.1496988574979475
$ws.Range("J5").Value = 0.1496988574979476
$ws.Range("M5").Value = 19.42400133333333
$ws.Range("N5").Value = 58.272004
$ws.Range("O5").Value = 0.04177372766745037
$ws.Range("P5").Value = 0.04177372766745036
$ws.Range("Q5").Value = 156.9543802139133
$ws.Range("R5").Value = 1412.58942192522
$ws.Range("S5").Value = 0.006253479305247721
$ws.Range("T5").Value = 0.006253479305247721
$ws.Range("G6").Value = 8.080435
$ws.Range("H6").Value = 24.241305
$ws.Range("I6").Value = 0.1496988574979475
$ws.Range("J6").Value = 0.1496988574979476
$ws.Range("M6").Value = 129.3233566666667
$ws.Range("N6").Value = 387.97007
$ws.Range("O6").Value = 0.2781259427306063
$ws.Range("P6").Value = 0.2781259427306062
$ws.Range("Q6").Value = 1044.988977526817
$ws.Range("R6").Value = 9404.900797741351
$ws.Range("S6").Value = 0.04163513586731134
$ws.Range("T6").Value = 0.04163513586731134
$ws.Range("I7").Value = 0.2404784903431001
$ws.Range("J7").Value = 0.2404784903431001
$ws.Range("M7").Value = 6.305846
$ws.Range("N7").Value = 18.917538
$ws.Range("O7").Value = 0.01356150511917599
$ws.Range("P7").Value = 0.01356150511917599
$ws.Range("Q7").Value = 81.85323579007201
$ws.Range("R7").Value = 736.6791221106481
$ws.Range("S7").Value = 0.003261250277839666
$ws.Range("T7").Value = 0.003261250277839666
$ws.Range("I8").Value = 0.2404784903431001
$ws.Range("J8").Value = 0.2404784903431001
$ws.Range("O8").Value = 0.392557056479861
$ws.Range("P8").Value = 0.3925570564798609
$ws.Range("S8").Value = 0.09440152831580807
$ws.Range("T8").Value = 0.09440152831580806
$ws.Range("I9").Value = 0.2404784903431001
$ws.Range("J9").Value = 0.2404784903431001
$ws.Range("M9").Value = 127.396393
$ws.Range("N9").Value = 382.189179
$ws.Range("O9").Value = 0.2739817680029065
$ws.Range("P9").Value = 0.2739817680029065
$ws.Range("Q9").Value = 1653.672956021076
$ws.Range("R9").Value = 14883.05660418969
$ws.Range("S9").Value = 0.06588672195087246
$ws.Range("T9").Value = 0.06588672195087246
$ws.Range("I10").Value = 0.2404784903431001
$ws.Range("J10").Value = 0.2404784903431001
$ws.Range("M10").Value = 19.42400133333333
$ws.Range("N10").Value = 58.272004
$ws.Range("O10").Value = 0.04177372766745037
$ws.Range("P10").Value = 0.04177372766745036
$ws.Range("Q10").Value = 252.133870875376
$ws.Range("R10").Value = 2269.204837878384
$ws.Range("S10").Value = 0.01004568296547226
$ws.Range("T10").Value = 0.01004568296547226
$ws.Range("I11").Value = 0.2404784903431001
$ws.Range("J11").Value = 0.2404784903431001
$ws.Range("M11").Value = 129.3233566666667
$ws.Range("N11").Value = 387.97007
$ws.Range("O11").Value = 0.2781259427306063
$ws.Range("P11").Value = 0.2781259427306062
$ws.Range("Q11").Value = 1678.68596955908
$ws.Range("R11").Value = 15108.17372603172
$ws.Range("S11").Value = 0.06688330683310773
$ws.Range("T11").Value = 0.06688330683310771
$ws.Range("G12").Value = 15.25749233333333
$ws.Range("H12").Value = 45.772477
$ws.Range("I12").Value = 0.2826616599952471
$ws.Range("J12").Value = 0.2826616599952471
$ws.Range("M12").Value = 6.305846
$ws.Range("N12").Value = 18.917538
$ws.Range("O12").Value = 0.01356150511917599
$ws.Range("P12").Value = 0.01356150511917599
$ws.Range("Q12").Value = 96.21139700018067
$ws.Range("R12").Value = 865.902573001626
$ws.Range("S12").Value = 0.003833317549020326
$ws.Range("T12").Value = 0.003833317549020325
$ws.Range("G13").Value = 15.25749233333333
$ws.Range("H13").Value = 45.772477
$ws.Range("I13").Value = 0.2826616599952471
$ws.Range("J13").Value = 0.2826616599952471
$ws.Range("O13").Value = 0.392557056479861
$ws.Range("P13").Value = 0.3925570564798609
$ws.Range("Q13").Value = 2784.975743791269
$ws.Range("R13").Value = 25064.78169412142
$ws.Range("S13").Value = 0.1109608292274455
$ws.Range("T13").Value = 0.1109608292274454
$ws.Range("G14").Value = 15.25749233333333
$ws.Range("H14").Value = 45.772477
$ws.Range("I14").Value = 0.2826616599952471
$ws.Range("J14").Value = 0.2826616599952471
$ws.Range("M14").Value = 127.396393
$ws.Range("N14").Value = 382.189179
$ws.Range("O14").Value = 0.2739817680029065
$ws.Range("P14").Value = 0.2739817680029065
$ws.Range("Q14").Value = 1943.74948949182
$ws.Range("R14").Value = 17493.74540542638
$ws.Range("S14").Value = 0.07744414135213422
$ws.Range("T14").Value = 0.07744414135213422
$ws.Range("G15").Value = 15.25749233333333
$ws.Range("H15").Value = 45.772477
$ws.Range("I15").Value = 0.2826616599952471
$ws.Range("J15").Value = 0.2826616599952471
$ws.Range("M15").Value = 19.42400133333333
$ws.Range("N15").Value = 58.272004
$ws.Range("O15").Value = 0.04177372766745037
$ws.Range("P15").Value = 0.04177372766745036
$ws.Range("Q15").Value = 296.3615514259898
$ws.Range("R15").Value = 2667.253962833908
$ws.Range("S15").Value = 0.0118078312066709
$ws.Range("T15").Value = 0.0118078312066709
$ws.Range("G16").Value = 15.25749233333333
$ws.Range("H16").Value = 45.772477
$ws.Range("I16").Value = 0.2826616599952471
$ws.Range("J16").Value = 0.2826616599952471
$ws.Range("M16").Value = 129.3233566666667
$ws.Range("N16").Value = 387.97007
$ws.Range("O16").Value = 0.2781259427306063
$ws.Range("P16").Value = 0.2781259427306062
$ws.Range("Q16").Value = 1973.150122862599
$ws.Range("R16").Value = 17758.35110576339
$ws.Range("S16").Value = 0.07861554065997618
$ws.Range("T16").Value = 0.07861554065997617
$ws.Range("G17").Value = 4.142925
$ws.Range("H17").Value = 12.428775
$ws.Range("I17").Value = 0.07675219702895753
$ws.Range("J17").Value = 0.07675219702895753
$ws.Range("M17").Value = 6.305846
$ws.Range("N17").Value = 18.917538
$ws.Range("O17").Value = 0.01356150511917599
$ws.Range("P17").Value = 0.01356150511917599
$ws.Range("Q17").Value = 26.12464703955
$ws.Range("R17").Value = 235.12182335595
$ws.Range("S17").Value = 0.001040875312916212
$ws.Range("T17").Value = 0.001040875312916212
$ws.Range("G18").Value = 4.142925
$ws.Range("H18").Value = 12.428775
$ws.Range("I18").Value = 0.07675219702895753
$ws.Range("J18").Value = 0.07675219702895753
$ws.Range("O18").Value = 0.392557056479861
$ws.Range("P18").Value = 0.3925570564798609
$ws.Range("Q18").Value = 756.215070030825
$ws.Range("R18").Value = 6805.935630277425
$ws.Range("S18").Value = 0.0301296165440499
$ws.Range("T18").Value = 0.03012961654404989
$ws.Range("G19").Value = 4.142925
$ws.Range("H19").Value = 12.428775
$ws.Range("I19").Value = 0.07675219702895753
$ws.Range("J19").Value = 0.07675219702895753
$ws.Range("M19").Value = 127.396393
$ws.Range("N19").Value = 382.189179
$ws.Range("O19").Value = 0.2739817680029065
$ws.Range("P19").Value = 0.2739817680029065
$ws.Range("Q19").Value = 527.793701469525
$ws.Range("R19").Value = 4750.143313225725
$ws.Range("S19").Value = 0.02102870264010121
$ws.Range("T19").Value = 0.02102870264010121
$ws.Range("G20").Value = 4.142925
$ws.Range("H20").Value = 12.428775
$ws.Range("I20").Value = 0.07675219702895753
$ws.Range("J20").Value = 0.07675219702895753
$ws.Range("M20").Value = 19.42400133333333
$ws.Range("N20").Value = 58.272004
$ws.Range("O20").Value = 0.04177372766745037
$ws.Range("P20").Value = 0.04177372766745036
$ws.Range("Q20").Value = 80.4721807239
$ws.Range("R20").Value = 724.2496265151
$ws.Range("S20").Value = 0.003206225376566165
$ws.Range("T20").Value = 0.003206225376566165
$ws.Range("G21").Value = 4.142925
$ws.Range("H21").Value = 12.428775
$ws.Range("I21").Value = 0.07675219702895753
$ws.Range("J21").Value = 0.07675219702895753
$ws.Range("M21").Value = 129.3233566666667
$ws.Range("N21").Value = 387.97007
$ws.Range("O21").Value = 0.2781259427306063
$ws.Range("P21").Value = 0.2781259427306062
$ws.Range("Q21").Value = 535.77696741825
$ws.Range("R21").Value = 4821.99270676425
$ws.Range("S21").Value = 0.02134677715532405
$ws.Range("T21").Value = 0.02134677715532404
$ws.Range("G22").Value = 13.51654933333334
$ws.Range("H22").Value = 40.549648
$ws.Range("I22").Value = 0.2504087951347477
$ws.Range("J22").Value = 0.2504087951347477
$ws.Range("M22").Value = 6.305846
$ws.Range("N22").Value = 18.917538
$ws.Range("O22").Value = 0.01356150511917599
$ws.Range("P22").Value = 0.01356150511917599
$ws.Range("Q22").Value = 85.23327854740268
$ws.Range("R22").Value = 767.0995069266241
$ws.Range("S22").Value = 0.003395920157106572
$ws.Range("T22").Value = 0.003395920157106571
$ws.Range("G23").Value = 13.51654933333334
$ws.Range("H23").Value = 40.549648
$ws.Range("I23").Value = 0.2504087951347477
$ws.Range("J23").Value = 0.2504087951347477
$ws.Range("O23").Value = 0.392557056479861
$ws.Range("P23").Value = 0.3925570564798609
$ws.Range("Q23").Value = 2467.198489154829
$ws.Range("R23").Value = 22204.78640239346
$ws.Range("S23").Value = 0.09829973953476509
$ws.Range("T23").Value = 0.09829973953476508
$ws.Range("G24").Value = 13.51654933333334
$ws.Range("H24").Value = 40.549648
$ws.Range("I24").Value = 0.2504087951347477
$ws.Range("J24").Value = 0.2504087951347477
$ws.Range("M24").Value = 127.396393
$ws.Range("N24").Value = 382.189179
$ws.Range("O24").Value = 0.2739817680029065
$ws.Range("P24").Value = 0.2739817680029065
$ws.Range("Q24").Value = 1721.959630873222
$ws.Range("R24").Value = 15497.63667785899
$ws.Range("S24").Value = 0.06860744441449579
$ws.Range("T24").Value = 0.06860744441449579
$ws.Range("G25").Value = 13.51654933333334
$ws.Range("H25").Value = 40.549648
$ws.Range("I25").Value = 0.2504087951347477
$ws.Range("J25").Value = 0.2504087951347477
$ws.Range("M25").Value = 19.42400133333333
$ws.Range("N25").Value = 58.272004
$ws.Range("O25").Value = 0.04177372766745037
$ws.Range("P25").Value = 0.04177372766745036
$ws.Range("Q25").Value = 262.5454722727325
$ws.Range("R25").Value = 2362.909250454592
$ws.Range("S25").Value = 0.01046050881349332
$ws.Range("T25").Value = 0.01046050881349332
$ws.Range("G26").Value = 13.51654933333334
$ws.Range("H26").Value = 40.549648
$ws.Range("I26").Value = 0.2504087951347477
$ws.Range("J26").Value = 0.2504087951347477
$ws.Range("M26").Value = 129.3233566666667
$ws.Range("N26").Value = 387.97007
$ws.Range("O26").Value = 0.2781259427306063
$ws.Range("P26").Value = 0.2781259427306062
$ws.Range("Q26").Value = 1748.005530337263
$ws.Range("R26").Value = 15732.04977303536
$ws.Range("S26").Value = 0.06964518221488694
$ws.Range("T26").Value = 0.06964518221488693

Write-Host "Updated 278 cells"
